$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item("2025-08-16")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate the most recent day sheet (keeps sheetPr/pageMargins/column widths/
# header style) then rename + overwrite with the new days ranking data.
$srcSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2025-08-17"

# Header row (unchanged text, kept for clarity)
$ws.Cells.Item(1,1).Value = "rank"
$ws.Cells.Item(1,2).Value = "title"
$ws.Cells.Item(1,3).Value = "author"
$ws.Cells.Item(1,4).Value = "latest_episode"

# New ranking data rows for 2025-08-17
$data = @(
  @(1, "いとこのこ", "いぬちく(著者)", "第37話"),
  @(2, "氷結令嬢さまをフォローしたら、メチャメチャ溺愛されてしまった件@comic", "漫画：ハレノチアメ 原作：愛坂タカト キャラクター原案：Bcoca", "第9話"),
  @(3, "リビルドワールド", "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)", "第72話②"),
  @(4, "宇崎ちゃんは遊びたい！", "丈(著者)", "第126話"),
  @(5, "サーシャちゃんとクラスメイトオタクくん", "はぐはぐ(著者)", "第83話"),
  @(6, "ダメ人間の愛しかた", "岩葉(著者)", "第19話前編　ダメ人間と新生活の彼女"),
  @(7, "魔都精兵のスレイブ", "原作:タカヒロ　漫画:竹村洋平", "第159話　神奴隷"),
  @(8, "聖液鍛冶屋のエロランタ", "しげきっくす(著者)", "第7話"),
  @(9, "みつばものがたり 呪いの少女と死の輪舞《ロンド》", "堤りん(漫画) 七沢またり(原作) EURA(キャラクター原案)", "第11話：勝利の美酒"),
  @(10, "えろいことするために巨乳美少女奴隷を買ったはずが、お師匠さまと慕われて思った通りにいかなくなる話", "佐藤36(作画) 煮豆シューター(原作)", "第3話後半"),
  @(11, "リアリスト魔王による聖域なき異世界改革", "鈴木マナツ(漫画) 羽田遼亮(原作) ゆーげん(キャラクターデザイン) ひたきゆう(キャラクターデザイン)", "第68幕②"),
  @(12, "転生してあらゆるモノに好かれながら異世界で好きな事をして生きて行く", "都尾琉(漫画) 御峰。(原作)", "第27話③"),
  @(13, "序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する", "作画：マエD 原作：新人", "第5話(4)"),
  @(14, "アラサーがVTuberになった話。", "犬威赤彦(漫画) とくめい(原作) カラスBTK(キャラクター原案)", "第26話"),
  @(15, "最強の少年聖騎士、転生者を狩る", "作画：御塩 原作：宇奈木ユラ", "第7話(1)"),
  @(16, "最強勇者パーティーは愛が知りたい", "山田肌襦袢", "第29話「きみがきみであればいい」"),
  @(17, "俺堕ちスレイブヒーローコレクション", "ゆっ栗栖(著者)", "第12話前半"),
  @(18, "役目を果たした日陰の勇者は、辺境で自由に生きていきます", "船野真帆(作画) 丘野優(原作) 布施龍太(キャラクター原案)", "第5話後半"),
  @(19, "愚かな天使は悪魔と踊る", "アズマサワヨシ(著者)", "第101話②"),
  @(20, "ハーレムより平穏を！異世界で静かにニート姫させてくれ", "さかたはるき(原作) かわやばぐ(作画)", "第14話後半"),
  @(21, "きみの願いが叶うまで", "浅月のりと(著者)", "第4話-2"),
  @(22, "よくわからないけれど異世界に転生していたようです", "内々けやき あし カオミン", "第137話 よくわからないけれど脱出するみたいです（２）"),
  @(23, "豚のレバーは加熱しろ", "みなみ(漫画) 逆井卓馬(原作) 遠坂あさぎ(キャラクターデザイン)", "第43話①"),
  @(24, "花子さん、食べないで", "茸谷きの子(著者)", "第9話後半"),
  @(25, "くらいあの子としたいこと", "碇マナツ(著者)", "第82話"),
  @(26, "小林さんちのメイドラゴン", "クール教信者", "第148話"),
  @(27, "追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。", "六志麻あさ 業務用餅 kisui", "第７０話"),
  @(28, "不徳のギルド", "河添太一", "第９７話：立派に育った所"),
  @(29, "治癒魔法の間違った使い方 ~戦場を駆ける回復要員~", "九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)", "第81話その3"),
  @(30, "十歳の最強魔導師", "猫月 天乃聖樹", "第2話（後編）"),
  @(31, "今日から僕は、彼女の✕✕を解消する。", "コアヤアコ(著者)", "第1話"),
  @(32, "異種族追放コンカフェ", "佐々木マサヒト(著者)", "第18話-1"),
  @(33, "王様ランキング200話～", "十日草輔（とおかそうすけ）", "第264話"),
  @(34, "理想のヒモ生活", "日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)", "第86話　その3"),
  @(35, "路地裏で拾った女の子がバッドエンド後の乙女ゲームのヒロインだった件", "カボチャマスク(原作) 樋乃えなが(作画) へいろー(キャラクター原案)", "第1話"),
  @(36, "おっさん、転生して天才役者になる", "芽々ノ圭(漫画) ほえ太郎(原作)", "第38話②"),
  @(37, "神への餌は彼女達", "蔵人幸明(原作) まりお(作画)", "第6話前半"),
  @(38, "お前妹じゃなくて許嫁だったのかよ!?", "湯猫子(漫画) 未来人A(原作)", "第29話"),
  @(39, "怪異部～M県Y市の怪現象について～", "さりい・Ｂ(著者)", "File.8"),
  @(40, "勇者パーティから追い出された不遇職【罠士】、ユニークスキル【矢印】で最強になる", "作画：たつひこ 原作：白石 有希", "第8話(3)"),
  @(41, "ゲーセン少女と異文化交流", "安原宏和(著者)", "第128話"),
  @(42, "終末ツーリング", "さいとー栄(著者)", "第49話 三沢基地　その５③"),
  @(43, "王都の外れの錬金術師 ～ハズレ職業だったので、のんびりお店経営します～", "あさなや(著者) yocco(原作) 純粋(キャラクター原案)", "element.50"),
  @(44, "うちの清楚系委員長がかつて中二病アイドルだったことを俺だけが知っている。", "三上こた こばやし少女 寝子空兄 ゆがー", "第1話　災禍の悪夢"),
  @(45, "辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？", "tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)", "第10話"),
  @(46, "ぽんドロイド！ はまさん", "はれやまはれぞう(著者)", "第6話"),
  @(47, "受付嬢に告白したくてギルドに通いつめたら英雄になってた", "よねぞう(作画) 猫又ぬこ(原作)", "第5話後半"),
  @(48, "不運からの最強男", "中林ずん(作画) フクフク(原作)", "第5話後半"),
  @(49, "農学博士の異世界無双～禁忌の知識で築くモンスター娘ハーレム～", "インド僧(原作) ヤスウミ(作画)", "第25話"),
  @(50, "最強で最速の無限レベルアップ ～スキル【経験値1000倍】と【レベルフリー】でレベル上限の枷が外れた俺は無双する～", "シオヤマ琴 鳥羽田 航 トモゼロ", "休載イラスト")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

$ws.Range("A1").Select()
Write-Output "Added sheet 2025-08-17 with $($data.Count) ranking rows"